# Updated cryptos list with refreshed prices / volume figures and a few reshuffled rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.474.48'
$ws.Range("E2").Value = '  +5.32%  '
$ws.Range("D3").Value = '3.331.66'
$ws.Range("E3").Value = '  +5.14%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'553.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.93%  '
$ws.Range("D6").Value = "'151.82"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.90%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = "'0.529"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.55%  '
$ws.Range("E9").Value = '  +3.79%  '
$ws.Range("E10").Value = '  +5.73%  '
$ws.Range("E11").Value = '  +2.13%  '
$ws.Range("D12").Value = '3.902.42'
$ws.Range("E12").Value = '  +5.10%  '
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("E14").Value = '  +4.59%  '
$ws.Range("E15").Value = '  +5.52%  '
$ws.Range("D16").Value = '62.392.86'
$ws.Range("E16").Value = '  +5.12%  '
$ws.Range("D17").Value = '3.322.01'
$ws.Range("E17").Value = '  +4.98%  '
$ws.Range("D18").Value = "'6.51"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +5.55%  '
$ws.Range("D19").Value = "'13.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.88%  '
$ws.Range("E20").Value = '  +5.05%  '
$ws.Range("D21").Value = "'384.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("D23").Value = "'0.537"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").Value = "'70.80"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.42%  '
$ws.Range("E25").Value = '  +3.62%  '
$ws.Range("D26").Value = "'8.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.27%  '
$ws.Range("D27").Value = '0.0₃0974'
$ws.Range("E27").Value = '  +9.88%  '
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("E29").Value = '  +3.77%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = "'6.35"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.64%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = "'1.32"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +11.76%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'22.95"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.22%  '
$ws.Range("E33").Value = '  +3.45%  '
$ws.Range("D34").Value = "'6.77"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.21%  '
$ws.Range("E35").Value = '  +11.47%  '
$ws.Range("D36").Value = "'159.34"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.96%  '
$ws.Range("D38").Value = "'26.91"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.85%  '
$ws.Range("D39").Value = '2.856.13'
$ws.Range("E39").Value = '  +5.04%  '
$ws.Range("D40").Value = "'0.0730"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.45%  '
$ws.Range("D41").Value = "'0.0317"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +9.04%  '
$ws.Range("E42").Value = '  +1.94%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = "'0.749"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.92%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = "'40.73"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.54%  '
$ws.Range("D45").Value = "'1.05"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.71%  '
$ws.Range("D46").Value = "'22.02"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +7.95%  '
$ws.Range("D47").Value = '3.375.39'
$ws.Range("E47").Value = '  +5.08%  '
$ws.Range("E48").Value = '  +4.93%  '
$ws.Range("D49").Value = "'6.30"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").Value = "'0.814"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +6.52%  '
$ws.Range("D51").Value = "'284.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +10.29%  '
